$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 826.86838530406374
$ws.Range("C2").Value = 461.50722525496786
$ws.Range("D2").Value = 975.62008339101214
$ws.Range("E2").Value = 417.24887536970186

$ws.Range("B3").Value = 770.17302114195684
$ws.Range("C3").Value = 444.87390981478114
$ws.Range("D3").Value = 1377.2977032292674
$ws.Range("E3").Value = 641.6741646289305

$ws.Range("B1:E3").Select() | Out-Null
